# Updates the cryptos list (Price / Volume(1h) columns) to the latest scraped
# values, matching a fresh run of the GitHub Actions scraper.
#
# Note: many "Price" values look like plain decimals (e.g. "1.003"), which
# Excel's COM layer would otherwise auto-convert to a floating point number
# (losing the original text formatting / trailing zeros). To keep them as
# literal text we prefix the assigned value with a leading apostrophe (the
# classic "treat as text" marker) and then reset the cell's Style back to
# "Normal" so the generated "quote prefixed" number format doesn't linger on
# the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''25.877.89'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '''1.638.76'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''215.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').Value = '''0.5057'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = '''1.003'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '''0.2579'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').Value = '''0.06437'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('D10').Value = '''19.73'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('D11').Value = '''0.07779'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '''4.292'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '''1.865.89'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '''1.637.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').Value = '''0.5614'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.03%  '
$ws.Range('D16').Value = '''0.0₅7632'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').Value = '''63.07'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').Value = '''25.894.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = '''1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = '''194.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').Value = '''4.328'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.38%  '
$ws.Range('D22').Value = '''9.912'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '''6.113'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').Value = '''1.003'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = '''1.775'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.56%  '
$ws.Range('D26').Value = '''140.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('D27').Value = '''0.1266'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').Value = '''6.837'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').Value = '''0.04881'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').Value = '''3.298'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.94%  '
$ws.Range('D33').Value = '''3.221'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').Value = '''1.569'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('D35').Value = '''2.376'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.45%  '
$ws.Range('D36').Value = '''0.9039'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').Value = '''2.579'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('D38').Value = '''0.5525'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('D39').Value = '''1.127.31'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').Value = '''0.01562'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('D41').Value = '''0.9969'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('D42').Value = '''5.549'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('D43').Value = '''0.8040'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = '''98.03'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').Value = '''1.776.80'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('D46').Value = '''0.0₈112'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.22%  '
$ws.Range('D47').Value = '''55.46'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('D48').Value = '''0.4278'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.07%  '
$ws.Range('D49').Value = '''7.715'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('D50').Value = '''0.05044'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.25%  '
$ws.Range('D51').Value = '''1.001'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.05%  '
